# Add team record (Wins/Losses/Ties) columns to the sheet, per commit:
# "Added team record to data" - W/L/T are added as columns on the same sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 1): AC=Wins, AD=Losses, AE=Ties
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the existing header formatting (bold, bordered, centered) by copying
# the format from an existing header cell onto the new header cells, rather
# than rebuilding the font/border/alignment by hand (which would mint a new
# style instead of reusing the workbook's existing header style).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data rows 2-38: Wins=63, Losses=99, Ties=0 for every team/player row.
$lastRow = 38
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 63  # AC
    $ws.Cells.Item($r, 30).Value = 99  # AD
    $ws.Cells.Item($r, 31).Value = 0   # AE
}
